# Applies the crypto price/volume refresh described in the commit
# "Updated cryptos list ... with GitHub Actions".
# Column D ("Price") holds numeric-looking text (e.g. "211.72") that must
# stay plain text -- Excel auto-coerces bare numeric strings to numbers on
# assignment, so we force text format first and restore the default style
# afterwards (matches the source cells, which carry no explicit style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.718.37'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.602.12'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.72'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.11%  '
$ws.Range('E6').Value = '  -0.52%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.01'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('E8').Value = '  +0.24%  '
$ws.Range('E9').Value = '  +0.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.74'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.44%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0842'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.64%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.826.80'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.595.40'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.25%  '
$ws.Range('E14').Value = '  +0.42%  '
$ws.Range('E15').Value = '  +0.57%  '
$ws.Range('E16').Value = '  +0.71%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.696.24'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.24%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0₃0746'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.58%  '
$ws.Range('E19').Value = '  +2.89%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '210.82'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.01%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.01'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.19%  '
$ws.Range('E22').Value = '  +1.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.31'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.41%  '
$ws.Range('E24').Value = '  +1.13%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.49'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.45%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('E27').Value = '  -0.46%  '
$ws.Range('E28').Value = '  -0.83%  '
$ws.Range('E29').Value = '  +1.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0515'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.57%  '
$ws.Range('E31').Value = '  -0.28%  '
$ws.Range('E32').Value = '  +1.67%  '
$ws.Range('E33').Value = '  +1.67%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.300.22'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.24%  '
$ws.Range('E35').Value = '  +0.50%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.611'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.27%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.16'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +22.51%  '
$ws.Range('E39').Value = '  -0.24%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.824'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.95%  '
$ws.Range('E41').Value = '  -1.49%  '
$ws.Range('E42').Value = '  -0.20%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.781'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.72%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '63.32'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.737.62'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '91.21'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.36%  '
$ws.Range('E47').Value = '  -2.28%  '
$ws.Range('E48').Value = '  -1.48%  '
$ws.Range('E49').Value = '  -0.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0519'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.06%  '
$ws.Range('E51').Value = '  +0.15%  '
